$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the data (excluding the header row) and sort it ascending by Year (col A).
$dataRange = $ws.Range("A2:C22")
$dataRange.Select() | Out-Null

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A22"), 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 0
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

$dataRange.Select() | Out-Null
